# Refresh the cryptocurrency price / 1h-volume table with the latest values.
# Rows 12/13 and 47/48/49 also got re-ranked (names + links swapped around).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.261.77"
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").Value = "'1.659.98"
$ws.Range("E3").Value = '  +0.54%  '

$ws.Range("E4").Value = '  +0.62%  '

$ws.Range("D5").Value = "'218.49"
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").Value = "'0.5335"
$ws.Range("E6").Value = '  +0.98%  '

$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("D8").Value = "'0.2635"
$ws.Range("E8").Value = '  +1.02%  '

$ws.Range("D9").Value = "'0.06359"

$ws.Range("D10").Value = "'20.55"
$ws.Range("E10").Value = '  +1.03%  '

$ws.Range("D11").Value = "'0.07841"
$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'1.706.75"
$ws.Range("E12").Value = '  +3.21%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.541"
$ws.Range("E13").Value = '  +1.56%  '

$ws.Range("D14").Value = "'1.886.81"
$ws.Range("E14").Value = '  +0.54%  '

$ws.Range("D15").Value = "'0.5518"
$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("D16").Value = "'0.0₅8194"
$ws.Range("E16").Value = '  +1.11%  '

$ws.Range("D17").Value = "'65.61"
$ws.Range("E17").Value = '  +0.87%  '

$ws.Range("D18").Value = "'26.237.32"
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").Value = "'1.009"
$ws.Range("E19").Value = '  +0.62%  '

$ws.Range("D20").Value = "'4.643"
$ws.Range("E20").Value = '  +2.08%  '

$ws.Range("D21").Value = "'192.13"
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").Value = "'10.13"
$ws.Range("E22").Value = '  +0.98%  '

$ws.Range("D23").Value = "'6.047"
$ws.Range("E23").Value = '  +1.04%  '

$ws.Range("D24").Value = "'1.010"
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("D25").Value = "'144.68"
$ws.Range("E25").Value = '  +3.30%  '

$ws.Range("D26").Value = "'0.1230"
$ws.Range("E26").Value = '  -0.77%  '

$ws.Range("D27").Value = "'7.231"
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  -0.58%  '

$ws.Range("E29").Value = '  +2.79%  '

$ws.Range("D30").Value = "'0.05798"
$ws.Range("E30").Value = '  -1.83%  '

$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").Value = "'3.576"
$ws.Range("E32").Value = '  +2.25%  '

$ws.Range("D33").Value = "'3.287"
$ws.Range("E33").Value = '  +1.52%  '

$ws.Range("D34").Value = "'1.610"
$ws.Range("E34").Value = '  +4.36%  '

$ws.Range("D35").Value = "'2.819"
$ws.Range("E35").Value = '  +2.48%  '

$ws.Range("D36").Value = "'0.9578"
$ws.Range("E36").Value = '  +1.80%  '

$ws.Range("D37").Value = "'2.430"
$ws.Range("E37").Value = '  +0.79%  '

$ws.Range("D38").Value = "'0.5790"
$ws.Range("E38").Value = '  +2.12%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").Value = "'5.846"
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("D41").Value = "'0.8532"
$ws.Range("E41").Value = '  +1.17%  '

$ws.Range("D42").Value = "'1.009"
$ws.Range("E42").Value = '  +0.63%  '

$ws.Range("D43").Value = "'104.70"
$ws.Range("E43").Value = '  +4.00%  '

$ws.Range("D44").Value = "'1.044.76"
$ws.Range("E44").Value = '  +3.74%  '

$ws.Range("D45").Value = "'1.799.87"
$ws.Range("E45").Value = '  +0.34%  '

$ws.Range("D46").Value = "'57.06"
$ws.Range("E46").Value = '  +0.39%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = "'0.0₈105"
$ws.Range("E47").Value = '  -1.36%  '

$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = "'1.011"
$ws.Range("E48").Value = '  +0.42%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = "'0.4370"
$ws.Range("E49").Value = '  +1.75%  '

$ws.Range("D50").Value = "'7.961"
$ws.Range("E50").Value = '  +1.71%  '

$ws.Range("D51").Value = "'0.05159"
$ws.Range("E51").Value = '  +0.17%  '
